$d = $word.ActiveDocument

function Set-ParagraphStrikeThrough($needle) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            $target = $p
            break
        }
    }
    if ($target -eq $null) {
        Write-Output "NOT FOUND: $needle"
        return
    }

    $pStart = $target.Range.Start
    $pEnd = $target.Range.End

    # Strike the whole paragraph (sets the paragraph-mark rPr + every plain run).
    $target.Range.Font.StrikeThrough = 1

    # Runs inside hyperlinks are not touched by the paragraph-range pass above,
    # so strike them explicitly too.
    foreach ($hl in $d.Hyperlinks) {
        $hlr = $hl.Range
        if ($hlr.Start -ge $pStart -and $hlr.End -le $pEnd) {
            $hlr.Font.StrikeThrough = 1
        }
    }
}

# "author: Alin Sinp | f04c0ce20 Execution,Interface,Message Routing: Fix bug in
#  StreamsTesting fixture (#1014)"
Set-ParagraphStrikeThrough("f04c0ce20 Execution,Interface,Message Routing: Fix bug in StreamsTesting fixture")

# "author: Leo Eich | 43c59b2ff Consensus,Interface: Make Cannot report master
#  public key changed metric warning less noisy (#986)"
Set-ParagraphStrikeThrough("43c59b2ff Consensus,Interface: Make Cannot report master public key changed metric warning less noisy")
